$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 261
$ws1.Range("F4").Value = 904
$ws1.Range("F6").Value = 44

# Sheet "全部类型" (All Types) - fourth sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 261
$ws4.Range("F5").Value = 904
$ws4.Range("F7").Value = 44
